$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.282.59'
$ws.Range('E2').Value = '  -1.14%  '
$ws.Range('D3').Value = '2.306.85'
$ws.Range('E3').Value = '  -2.12%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''317.49'
$ws.Range('E5').Value = '  -0.39%  '
$ws.Range('E6').Value = '  -1.91%  '
$ws.Range('D7').Value = '''0.629'
$ws.Range('E7').Value = '  -1.29%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = '''0.612'
$ws.Range('E9').Value = '  -1.62%  '
$ws.Range('E10').Value = '  -4.11%  '
$ws.Range('D11').Value = '''0.0910'
$ws.Range('E11').Value = '  -1.96%  '
$ws.Range('D12').Value = '''8.46'
$ws.Range('E12').Value = '  -0.09%  '
$ws.Range('E13').Value = '  +0.50%  '
$ws.Range('E14').Value = '  -1.92%  '
$ws.Range('D15').Value = '''15.50'
$ws.Range('E15').Value = '  -3.32%  '
$ws.Range('D16').Value = '2.655.43'
$ws.Range('E16').Value = '  -2.14%  '
$ws.Range('D17').Value = '2.301.18'
$ws.Range('E17').Value = '  -1.90%  '
$ws.Range('D18').Value = '42.214.96'
$ws.Range('E18').Value = '  -1.15%  '
$ws.Range('D19').Value = '''7.79'
$ws.Range('E19').Value = '  -0.05%  '
$ws.Range('E20').Value = '  -0.63%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = '''286.34'
$ws.Range('E21').Value = '  +11.00%  '
$ws.Range('B22').Value = 'Litecoin'
$ws.Range('C22').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D22').Value = '''73.93'
$ws.Range('E22').Value = '  -3.77%  '
$ws.Range('D23').Value = '''3.55'
$ws.Range('E23').Value = '  -1.42%  '
$ws.Range('E24').Value = '  -0.95%  '
$ws.Range('D25').Value = '''10.04'
$ws.Range('E25').Value = '  +5.66%  '
$ws.Range('E26').Value = '  +0.55%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').Value = '''11.00'
$ws.Range('E27').Value = '  -3.79%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '''23.46'
$ws.Range('E28').Value = '  +1.66%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '''2.23'
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').Value = '''35.63'
$ws.Range('E30').Value = '  -2.45%  '
$ws.Range('D31').Value = '''165.24'
$ws.Range('E31').Value = '  -5.58%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '''0.0884'
$ws.Range('E32').Value = '  -0.64%  '
$ws.Range('B33').Value = 'WEMIXToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D33').Value = '''2.92'
$ws.Range('E33').Value = '  -1.08%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '''5.94'
$ws.Range('E34').Value = '  -3.80%  '
$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D35').Value = '''0.133'
$ws.Range('E35').Value = '  +1.34%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').Value = '''0.119'
$ws.Range('E36').Value = '  -6.95%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').Value = '''4.68'
$ws.Range('E37').Value = '  +1.04%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').Value = '''2.95'
$ws.Range('E38').Value = '  +10.00%  '
$ws.Range('D39').Value = '''0.0354'
$ws.Range('E39').Value = '  -2.26%  '
$ws.Range('B40').Value = 'NEARProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D40').Value = '''3.65'
$ws.Range('E40').Value = '  -4.10%  '
$ws.Range('B41').Value = 'BitcoinSV'
$ws.Range('C41').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D41').Value = '''102.34'
$ws.Range('E41').Value = '  +19.48%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').Value = '''1.50'
$ws.Range('E42').Value = '  +1.80%  '
$ws.Range('B43').Value = 'MultiversX'
$ws.Range('C43').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D43').Value = '''71.36'
$ws.Range('E43').Value = '  -0.76%  '
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').Value = '''0.228'
$ws.Range('E44').Value = '  -5.15%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').Value = '''1.00'
$ws.Range('E45').Value = '  +0.30%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '''116.90'
$ws.Range('E46').Value = '  +2.32%  '
$ws.Range('B47').Value = 'Celestia'
$ws.Range('C47').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D47').Value = '''12.20'
$ws.Range('E47').Value = '  +1.22%  '
$ws.Range('B48').Value = 'ordi'
$ws.Range('C48').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D48').Value = '''78.88'
$ws.Range('E48').Value = '  +4.14%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').Value = '''9.19'
$ws.Range('E49').Value = '  +0.51%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').Value = '''5.36'
$ws.Range('E50').Value = '  -2.52%  '
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').Value = '''1.29'
$ws.Range('E51').Value = '  +2.09%  '
